{"js": "// Remove the pre-existing \"_GoBack\" bookmark (it sat, empty, at the end of\n// the \"-Fix the 3.3 label to 5V...\" paragraph). It gets re-created below,\n// now spanning the \"MSP 430 To do list\" heading instead.\ncontext.document.deleteBookmark(\"_GoBack\");\nawait context.sync();\n\n// Drop the three stale checklist items that were finished/removed from the\n// list (\"Incorporate ferrite bead footprint\", \"Resitor footprint\",\n// \"Capacitor Footprint\"). They were paragraphs 2-4 (1-based), i.e. indexes\n// 1-3 right after the title paragraph.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nparagraphs.items[1].delete();\nparagraphs.items[2].delete();\nparagraphs.items[3].delete();\nawait context.sync();\n\n// Underline the title paragraph (\"MSP 430 To do list\") - both the runs and\n// the paragraph mark itself - and re-insert the \"_GoBack\" bookmark so it\n// spans the whole heading.\nconst titleParagraphs = body.paragraphs;\ntitleParagraphs.load(\"items\");\nawait context.sync();\n\nconst title = titleParagraphs.items[0];\ntitle.font.underline = Word.UnderlineType.single;\n\nconst titleRange = title.getRange(\"Whole\");\ntitleRange.insertBookmark(\"_GoBack\");\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Remove the pre-existing \"_GoBack\" bookmark (it sat, empty, at the end of\n# the \"-Fix the 3.3 label to 5V...\" paragraph). It gets re-created below,\n# now spanning the \"MSP 430 To do list\" heading instead.\nif ($d.Bookmarks.Exists(\"_GoBack\")) {\n    $d.Bookmarks.Item(\"_GoBack\").Delete()\n}\n\n# Drop the three stale checklist items that were finished/removed from the\n# list (\"Incorporate ferrite bead footprint\", \"Resitor footprint\",\n# \"Capacitor Footprint\"). They are paragraphs 2-4 (1-based), right after the\n# title paragraph.\n$d.Paragraphs.Item(2).Range.Delete()\n$d.Paragraphs.Item(2).Range.Delete()\n$d.Paragraphs.Item(2).Range.Delete()\n\n# Underline the title paragraph (\"MSP 430 To do list\") - both the runs and\n# the paragraph mark itself - and re-insert the \"_GoBack\" bookmark so it\n# spans the whole heading.\n$title = $d.Paragraphs.Item(1).Range\n$title.Font.Underline = 1  # wdUnderlineSingle\n\n$d.Bookmarks.Add(\"_GoBack\", $title)\n"}
